$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels for D1 and E1 (mean/sd -> log10 mean/log10 sd)
$ws.Range("D1").Value = "log10 mean"
$ws.Range("E1").Value = "log10 sd"

# Replace mean/sd data values (D2:E37) with their log10-transformed equivalents
$ws.Range("D2").Value = 9.0815117409882
$ws.Range("E2").Value = 0.560883533000909
$ws.Range("D3").Value = 8.83388194420938
$ws.Range("E3").Value = 0.267541478270403
$ws.Range("D4").Value = 8.98597995386476
$ws.Range("E4").Value = 0.290428882522969
$ws.Range("D5").Value = 8.98872751732753
$ws.Range("E5").Value = 0.203236869456351
$ws.Range("D6").Value = 8.97387684729539
$ws.Range("E6").Value = 0.252040053510732
$ws.Range("D7").Value = 8.54526150028628
$ws.Range("E7").Value = 0.839993281753799
$ws.Range("D8").Value = 9.36773922902156
$ws.Range("E8").Value = 0.217127470661799
$ws.Range("D9").Value = 8.64600229814254
$ws.Range("E9").Value = 0.39391540053294
$ws.Range("D10").Value = 8.64389048112528
$ws.Range("E10").Value = 0.305681305317148
$ws.Range("D11").Value = 8.79635927508536
$ws.Range("E11").Value = 0.18842981107003
$ws.Range("D12").Value = 8.62183536283246
$ws.Range("E12").Value = 0.548299732979768
$ws.Range("D13").Value = 9.24763830120195
$ws.Range("E13").Value = 0.324261188151481
$ws.Range("D14").Value = 8.95708935459208
$ws.Range("E14").Value = 0.581938428066334
$ws.Range("D15").Value = 8.58703947426276
$ws.Range("E15").Value = 0.338464209818941
$ws.Range("D16").Value = 8.70992662630511
$ws.Range("E16").Value = 0.591838602804627
$ws.Range("D17").Value = 8.56853077451129
$ws.Range("E17").Value = 0.482960529406324
$ws.Range("D18").Value = 8.79562276815422
$ws.Range("E18").Value = 0.363676289854048
$ws.Range("D19").Value = 8.73629920907621
$ws.Range("E19").Value = 0.440651556198468
$ws.Range("D20").Value = 8.86669431343842
$ws.Range("E20").Value = 0.240298328774395
$ws.Range("D21").Value = 8.58221371418316
$ws.Range("E21").Value = 0.522512154265031
$ws.Range("D22").Value = 8.31352584981283
$ws.Range("E22").Value = 1.09456374857094
$ws.Range("D23").Value = 8.61061766468389
$ws.Range("E23").Value = 0.35290039597624
$ws.Range("D24").Value = 7.87040232763843
$ws.Range("E24").Value = 0.315397000081903
$ws.Range("D25").Value = 8.67759207997035
$ws.Range("E25").Value = 0.203368238740061
$ws.Range("D26").Value = 7.851104188672
$ws.Range("E26").Value = 0.828796074880238
$ws.Range("D27").Value = 8.25511230711757
$ws.Range("E27").Value = 0.25442704315756
$ws.Range("D28").Value = 6.96154005980229
$ws.Range("E28").Value = 0.701148845461913
$ws.Range("D29").Value = 6.32417117848861
$ws.Range("E29").Value = 1.38762828456548
$ws.Range("D30").Value = 3.69141973923717
$ws.Range("E30").Value = 2.37455044756191
$ws.Range("D31").Value = 5.89871449281301
$ws.Range("E31").Value = 0.838868261159922
$ws.Range("D32").Value = 6.5532249182237
$ws.Range("E32").Value = 0.967169048206081
$ws.Range("D33").Value = 6.36524959517012
$ws.Range("E33").Value = 1.21720612288495
$ws.Range("D34").Value = 6.13894138767063
$ws.Range("E34").Value = 1.49150826315435
$ws.Range("D35").Value = 7.86550424222228
$ws.Range("E35").Value = 0.438440659900706
$ws.Range("D36").Value = 0
$ws.Range("D37").Value = 0
